$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-32 down to 4-33
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly price-report entry
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44503
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112032
$ws.Range("G3").Value = "Zapallo italiano"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6750
$ws.Range("N3").Value = "$/caja 60 unidades"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 112
$ws.Range("Q3").Value = 60
$ws.Range("R3").Value = "Hortaliza"
